$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The scrape that produced this sheet had rows 2 and 4 out of order (the
# 22:40 GRU->NAT flight was listed before the 12:30 GRU->NAT flight) and the
# "Preco" (price) / "Moeda" (currency) columns (H/I) were left blank for
# every flight. This edit:
#   1) puts the 12:30 flight back in row 2 and the 22:40 flight in row 4
#      (row 4 keeps the taller, wrapped "2:00 / +1" arrival-time look that
#      row 2 used to have),
#   2) fills in the missing price/currency values for all four flights.
# ---------------------------------------------------------------------------

# Row 2 currently carries a wrap-text-derived row height (from the old "22:40"
# flight's multi-line arrival time). Delete + re-insert it so the fresh row
# starts out with the sheet's normal (non custom) row height before we put
# the 12:30 flight's data into it.
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "qua., 15/11"
$ws.Range("B2").Value = "12:30"
$ws.Range("C2").Value = "GRU"
$ws.Range("D2").Value = "15:50"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").WrapText = $true
$ws.Range("E2").Value = "NAT"
$ws.Range("F2").Value = "3 h 20 min."
$ws.Range("G2").Value = "Direto"

# Row 4 now becomes the 22:40 flight with the wrapped two-line arrival time,
# so it picks up the taller row height.
$ws.Range("B4").Value = "22:40"
$ws.Range("D4").Value = "2:00" + [char]13 + [char]10 + "+1"
$ws.Rows.Item(4).RowHeight = 28.8

# Fill in Preco (H) / Moeda (I) for every flight.
$ws.Range("H2").Value = "72.091"
$ws.Range("I2").Value = "pontos"

$ws.Range("H3").Value = "98.876"
$ws.Range("I3").Value = "pontos"

$ws.Range("H4").Value = "72.091"
$ws.Range("I4").Value = "pontos"

$ws.Range("H5").Value = "80.105"
$ws.Range("I5").Value = "pontos"

# Column H ("Preco") needs to be a touch wider to fit the new values.
$ws.Columns.Item(8).ColumnWidth = 5.666666666666667
